# B6-PowerPoint.pptx edit: 11 Jul 2020
#
# 1) Re-style the three summary tables (slides 14-16) from the custom
#    "Table_0" style to the built-in "No Style, Table Grid" style.
# 2) Re-colour the deck's theme from the "Integral" (Red Violet) palette
#    to the standard "Office Theme" palette.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Tables: switch from the custom style to the built-in
#    "No Style, Table Grid" style ({C59F6BC2-A839-4A0D-9F63-C5825AE57DDA})
# ---------------------------------------------------------------------
$noStyleTableGrid = "{C59F6BC2-A839-4A0D-9F63-C5825AE57DDA}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($noStyleTableGrid)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Theme: recolour to the Office Theme scheme.
# ---------------------------------------------------------------------
function ToBGRInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# MsoThemeColorSchemeIndex order: Dark1, Light1, Dark2, Light2,
# Accent1-6, Hyperlink, FollowedHyperlink.
$officeThemeColors = @(
    "000000", # Dark 1
    "FFFFFF", # Light 1
    "44546A", # Dark 2
    "E7E6E6", # Light 2
    "5B9BD5", # Accent 1
    "ED7D31", # Accent 2
    "A5A5A5", # Accent 3
    "FFC000", # Accent 4
    "4472C4", # Accent 5
    "70AD47", # Accent 6
    "0563C1", # Hyperlink
    "954F72"  # Followed Hyperlink
)

$theme = $p.SlideMaster.Theme
try { $theme.Name = "Office Theme" } catch { }

$themeColorScheme = $theme.ThemeColorScheme
for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $themeColorScheme.Item($i + 1).RGB = ToBGRInt $officeThemeColors[$i]
}
